# Trade #50 closed at 2026-02-17 15:42:34 - unknown UNKNOWN +0.000%

$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.42   # Current Capital
$summary.Range("B4").Value = 0.42      # Total P&L $
$summary.Range("B5").Value = 0.17      # Total P&L %
$summary.Range("B6").Value = 50        # Total Trades
$summary.Range("B8").Value = 27        # Losing Trades
$summary.Range("B9").Value = 28        # Win Rate %

# --- Strategy Status sheet (MarketMaking row) ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.42     # Capital
$status.Range("D4").Value = 50         # Trades
$status.Range("E4").Value = 0.42       # P&L $
$status.Range("F4").Value = 0.42       # P&L %
$status.Range("G4").Value = 28         # Win Rate %

# --- Trade #50 row (row 51) was closed; update on both "All Trades" and "MarketMaking" sheets ---
foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("G51").Value = 0.44
    $ws.Range("H51").Value = "CLOSED"
    $ws.Range("I51").Value = -18.5185
    $ws.Range("J51").Value = -0.1
    $ws.Range("K51").Value = 100.42
    $ws.Range("P51").Value = "early_exit"
    $ws.Range("Q51").Value = 6.61
}
